$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "28.943.58"
$ws.Range("E2").Value = "  -0.31%  "

# Row 3
$ws.Range("D3").Value = "1.828.72"
$ws.Range("E3").Value = "  -0.10%  "

# Row 4
$ws.Range("D4").Value = "'0.9968"
$ws.Range("E4").Value = "  -0.26%  "

# Row 5
$ws.Range("D5").Value = "'241.61"
$ws.Range("E5").Value = "  +0.11%  "

# Row 6
$ws.Range("D6").Value = "'0.6260"
$ws.Range("E6").Value = "  -4.14%  "

# Row 7
$ws.Range("D7").Value = "'0.9964"
$ws.Range("E7").Value = "  -0.39%  "

# Row 8
$ws.Range("D8").Value = "'0.07591"
$ws.Range("E8").Value = "  +3.53%  "

# Row 9
$ws.Range("E9").Value = "  -0.49%  "

# Row 10
$ws.Range("D10").Value = "'22.54"
$ws.Range("E10").Value = "  -1.72%  "

# Row 11
$ws.Range("D11").Value = "'0.07699"
$ws.Range("E11").Value = "  +0.36%  "

# Row 12
$ws.Range("D12").Value = "1.837.96"
$ws.Range("E12").Value = "  +0.25%  "

# Row 13
$ws.Range("D13").Value = "'4.942"
$ws.Range("E13").Value = "  -0.69%  "

# Row 14
$ws.Range("D14").Value = "'0.6634"
$ws.Range("E14").Value = "  -0.42%  "

# Row 15
$ws.Range("D15").Value = "'0.00001024"
$ws.Range("E15").Value = "  +17.82%  "

# Row 16
$ws.Range("D16").Value = "'82.62"
$ws.Range("E16").Value = "  +1.14%  "

# Row 17
$ws.Range("D17").Value = "'6.047"
$ws.Range("E17").Value = "  -0.82%  "

# Row 18
$ws.Range("D18").Value = "28.982.29"
$ws.Range("E18").Value = "  -0.24%  "

# Row 19
$ws.Range("D19").Value = "'225.70"
$ws.Range("E19").Value = "  +0.89%  "

# Row 20
$ws.Range("D20").Value = "'12.32"

# Row 21
$ws.Range("D21").Value = "'0.9955"
$ws.Range("E21").Value = "  -0.45%  "

# Row 22
$ws.Range("D22").Value = "'7.180"
$ws.Range("E22").Value = "  +1.04%  "

# Row 23
$ws.Range("D23").Value = "'0.9964"
$ws.Range("E23").Value = "  -0.35%  "

# Row 24
$ws.Range("D24").Value = "'158.05"
$ws.Range("E24").Value = "  +0.06%  "

# Row 25
$ws.Range("D25").Value = "'8.481"
$ws.Range("E25").Value = "  -0.16%  "

# Row 26
$ws.Range("D26").Value = "'0.1368"
$ws.Range("E26").Value = "  -0.49%  "

# Row 27
$ws.Range("D27").Value = "'17.86"
$ws.Range("E27").Value = "  -0.19%  "

# Row 28
$ws.Range("D28").Value = "'1.482"
$ws.Range("E28").Value = "  -1.68%  "

# Row 29
$ws.Range("D29").Value = "'4.090"
$ws.Range("E29").Value = "  -0.36%  "

# Row 30
$ws.Range("D30").Value = "'4.008"
$ws.Range("E30").Value = "  -0.13%  "

# Row 31
$ws.Range("D31").Value = "'1.186"
$ws.Range("E31").Value = "  -1.12%  "

# Row 32
$ws.Range("D32").Value = "'0.05219"
$ws.Range("E32").Value = "  -2.39%  "

# Row 33
$ws.Range("D33").Value = "'1.837"
$ws.Range("E33").Value = "  +0.03%  "

# Row 34
$ws.Range("D34").Value = "'0.7346"
$ws.Range("E34").Value = "  -1.13%  "

# Row 35
$ws.Range("E35").Value = "  -1.59%  "

# Row 36
$ws.Range("D36").Value = "'2.683"
$ws.Range("E36").Value = "  +1.48%  "

# Row 37
$ws.Range("D37").Value = "1.235.33"
$ws.Range("E37").Value = "  -4.67%  "

# Row 39
$ws.Range("D39").Value = "'0.01783"
$ws.Range("E39").Value = "  -0.32%  "

# Row 40
$ws.Range("D40").Value = "'6.346"
$ws.Range("E40").Value = "  -0.19%  "

# Row 41
$ws.Range("D41").Value = "'0.8937"
$ws.Range("E41").Value = "  -0.40%  "

# Row 42
$ws.Range("D42").Value = "'0.9968"
$ws.Range("E42").Value = "  -0.25%  "

# Row 43
$ws.Range("D43").Value = "'101.62"
$ws.Range("E43").Value = "  -1.55%  "

# Row 44
$ws.Range("B44").Value = "BabyDogeCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D44").Value = "'0.00000000125"
$ws.Range("E44").Value = "  +4.03%  "

# Row 45
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").Value = "1.980.97"
$ws.Range("E45").Value = "  -0.19%  "

# Row 46
$ws.Range("D46").Value = "'63.95"
$ws.Range("E46").Value = "  -0.39%  "

# Row 47
$ws.Range("D47").Value = "'0.5087"
$ws.Range("E47").Value = "  -1.01%  "

# Row 48
$ws.Range("D48").Value = "'0.4032"
$ws.Range("E48").Value = "  +1.21%  "

# Row 49
$ws.Range("D49").Value = "'8.848"
$ws.Range("E49").Value = "  +1.54%  "

# Row 50
$ws.Range("D50").Value = "'0.05748"
$ws.Range("E50").Value = "  -1.43%  "

# Row 51
$ws.Range("D51").Value = "'6.681"
$ws.Range("E51").Value = "  -0.37%  "
